# Session 1: build derived bars_1m for ES (FULL + RTH)
# - Add ES_BARS_1M row to DATASETS sheet
# - Update ES row in INSTRUMENTS sheet (prices_dataset_id -> ES_BARS_1M, add OHLCV column mapping)

$wb = $excel.ActiveWorkbook

# --- DATASETS sheet: append new row 8 (row 7 left blank) ---
$dsSheet = $wb.Worksheets.Item("DATASETS")

$dsSheet.Range("A8").Value = "ES_BARS_1M"
$dsSheet.Range("B8").Value = "derived_bars"
$dsSheet.Range("C8").Value = "derived"
$dsSheet.Range("D8").Value = "DB_ES_OHLCV_1S"
$dsSheet.Range("E8").Value = "incremental"
$dsSheet.Range("F8").Value = "bar_time"
$dsSheet.Range("G8").Value = "UTC"
$dsSheet.Range("H8").Value = "1m"
$dsSheet.Range("I8").Value = "event_time"
$dsSheet.Range("J8").Value = 0
$dsSheet.Range("K8").Value = "bars_1m"
$dsSheet.Range("L8").Value = "instrument_id,session,date"
$dsSheet.Range("M8").Value = "instrument_id: ES"

# --- INSTRUMENTS sheet: update ES row (row 26) ---
$instSheet = $wb.Worksheets.Item("INSTRUMENTS")

$instSheet.Range("D26").Value = "ES_BARS_1M"
$instSheet.Range("E26").Value = "open"
$instSheet.Range("F26").Value = "high"
$instSheet.Range("G26").Value = "low"
$instSheet.Range("H26").Value = "close"
$instSheet.Range("Q26").Value = "volume"
